$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns H:I mirror the existing Exclude/Value pair pattern (e.g. F:G).
# Copy the header formatting (bold Calibri on pink fill) from F1:G1 onto H1:I1.
$ws.Range("F1:G1").Copy()
$ws.Range("H1:I1").PasteSpecial(-4122)

# Copy the data-row formatting from F2:G2 onto H2:I2.
$ws.Range("F2:G2").Copy()
$ws.Range("H2:I2").PasteSpecial(-4122)

$ws.Range("H1").Value = "Exclude4"
$ws.Range("I1").Value = "Value4"
$ws.Range("H2").Value = "brand_name"
$ws.Range("I2").Value = "RIGHT GUARD"

[void]$ws.Range("B2").Select()

Write-Output "done"
